# Adding new test scripts for watch list (TestCase_E33, TestCase_E34)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Append new test-case rows (55-58) ---
# Column C of row 58 is entered before column B/A of that row, matching the
# order the new shared strings were introduced upstream.

# Row 55
$ws.Cells.Item(55, 1).Value = "TestCase_B54"
$ws.Cells.Item(55, 2).Value = "OPQA-316"
$ws.Cells.Item(55, 3).Value = "Verify that only articles get displayed in the summary page when user searches using ARTICLES content type in search drop down"
$ws.Cells.Item(55, 4).Value = "Y"
$ws.Cells.Item(55, 5).Value = "SKIP"

# Row 56
$ws.Cells.Item(56, 1).Value = "TestCase_B55"
$ws.Cells.Item(56, 2).Value = "OPQA-366"
$ws.Cells.Item(56, 3).Value = "Verify that only patents get displayed in the summary page when user searches using PATENTS content type in search drop down"
$ws.Cells.Item(56, 4).Value = "Y"
$ws.Cells.Item(56, 5).Value = "SKIP"

# Row 57
$ws.Cells.Item(57, 1).Value = "TestCase_B56"
$ws.Cells.Item(57, 2).Value = "OPQA-744"
$ws.Cells.Item(57, 3).Value = "Verify that only posts get displayed in the summary page when user searches using POSTS content type in search drop down"
$ws.Cells.Item(57, 4).Value = "Y"
$ws.Cells.Item(57, 5).Value = "SKIP"

# Row 58 (C, then B, then A - matches the shared-string insertion order upstream)
$ws.Cells.Item(58, 3).Value = "Verify that only people get displayed in the summary page when user searches using people content type in search drop down"
$ws.Cells.Item(58, 2).Value = "OPQA-380"
$ws.Cells.Item(58, 1).Value = "TestCase_B57"
$ws.Cells.Item(58, 4).Value = "Y"
$ws.Cells.Item(58, 5).Value = "PASS"

# --- Formatting to match the rest of the table (thin borders all around,
#     wrapped description column) ---
$newRange = $ws.Range("A55:E58")
$newRange.Borders.LineStyle = 1
$ws.Range("C55:C58").WrapText = $true

# --- Sheet view: selection moves to D26, no more scrolled topLeftCell ---
$ws.Activate()
$ws.Range("D26").Select()
